# BusPirate v3.9 SSOP BOM edits
# - Resistor array RN1..RN3 -> RN1..RN4 (qty 3 -> 4)
# - IO connector CON1 (was "JST 1x10") -> 0.1" right-angled shrouded male header
# - USB/IO connector CON2 (was "JST 1x06") -> 0.1" straight male header
# - LED package LED-0603 -> LED-0804 (both LED rows)
# - Selection moved to D11
# - Column C widened to fit the new (longer) text

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: 10K 4 Resistor array / RN1, RN2, RN3 -> RN1, RN2, RN3, RN4 ; qty 3 -> 4
$ws.Range("E5").Value = 4
$ws.Range("F5").Value = "RN1, RN2, RN3, RN4"

# Row 18 (CON1): JST 1x10 connector replaced with 0.1" shrouded right-angled male header
$ws.Range("D18").Value = "1x10"
# Row 19 (CON2): JST 1x06 connector replaced with 0.1" male header
$ws.Range("D19").Value = "1x06"
$ws.Range("C19").Value = "0.1`" male header"
$ws.Range("C18").Value = "0.1`" shoruded right angled male header Molex 70553-0044"

# Rows 9-10: LED package LED-0603 -> LED-0804
$ws.Range("D9").Value = "LED-0804"
$ws.Range("D10").Value = "LED-0804"

# Widen column C to fit the new, longer connector description text
$ws.Columns.Item(3).ColumnWidth = 55.59

# Update the active selection cell, as recorded in the saved view
[void]$ws.Range("D11").Select()
